# Auto-generated Excel COM-interop script applying the Bahamut_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific
# Leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 34000
$ws.Range("I13").Value = 34000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 34000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -33831
$ws.Range("N13").ClearContents()

# Row 88
$ws.Range("H88").Value = 1224147.8
$ws.Range("I88").Value = 2343.1428
$ws.Range("J88").Value = 2649586.2
$ws.Range("K88").Value = 2343.1428
$ws.Range("L88").Value = 2649586.2
$ws.Range("M88").Value = -1937.1428
$ws.Range("N88").Value = -2650398.2

# Row 91
$ws.Range("H91").Value = 1224147.8
$ws.Range("I91").Value = 2343.1428
$ws.Range("J91").Value = 2649586.2
$ws.Range("K91").Value = 2343.1428
$ws.Range("L91").Value = 2649586.2
$ws.Range("M91").Value = -939.1428000000001
$ws.Range("N91").Value = -2652394.2

# Row 129
$ws.Range("H129").Value = 1278264.4
$ws.Range("I129").Value = 420
$ws.Range("J129").Value = 1611615
$ws.Range("K129").Value = 1260
$ws.Range("L129").Value = 4834845
$ws.Range("M129").Value = 3740
$ws.Range("N129").Value = -4844845

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7355.817
$ws.Range("I32").Value = 5425.224
$ws.Range("J32").Value = 15969.23
$ws.Range("K32").Value = 5425.224
$ws.Range("L32").Value = 15969.23
$ws.Range("M32").Value = -5138.224
$ws.Range("N32").Value = -16543.23

# Row 122
$ws.Range("H122").Value = 1553
$ws.Range("I122").Value = 1556
$ws.Range("J122").Value = 1550
$ws.Range("K122").Value = 4668
$ws.Range("L122").Value = 4650
$ws.Range("M122").Value = -2218
$ws.Range("N122").Value = -9550

$ws = $wb.Worksheets.Item("BSM")
# Row 14
$ws.Range("H14").Value = 55006
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 55006
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 55006
$ws.Range("N14").Value = -55350

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1060.1333
$ws.Range("I22").Value = 1222.5
$ws.Range("J22").Value = 410.66666
$ws.Range("K22").Value = 1222.5
$ws.Range("L22").Value = 410.66666
$ws.Range("M22").Value = -872.5
$ws.Range("N22").Value = -1110.66666

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

# Row 31
$ws.Range("H31").Value = 3060.1428
$ws.Range("I31").Value = 3538.4167
$ws.Range("J31").Value = 2016.6364
$ws.Range("K31").Value = 3538.4167
$ws.Range("L31").Value = 2016.6364
$ws.Range("M31").Value = -3243.4167
$ws.Range("N31").Value = -2606.6364

# Row 34
$ws.Range("H34").Value = 3060.1428
$ws.Range("I34").Value = 3538.4167
$ws.Range("J34").Value = 2016.6364
$ws.Range("K34").Value = 3538.4167
$ws.Range("L34").Value = 2016.6364
$ws.Range("M34").Value = -3336.4167
$ws.Range("N34").Value = -2420.6364

# Row 58
$ws.Range("H58").Value = 5390.9565
$ws.Range("I58").Value = 941.8823
$ws.Range("J58").Value = 17996.666
$ws.Range("K58").Value = 941.8823
$ws.Range("L58").Value = 17996.666
$ws.Range("M58").Value = -738.8823

# Row 62
$ws.Range("H62").Value = 6902.75
$ws.Range("I62").Value = 6985
$ws.Range("J62").Value = 6656
$ws.Range("K62").Value = 6985
$ws.Range("L62").Value = 6656
$ws.Range("M62").Value = -6361
$ws.Range("N62").Value = -7904

# Row 65
$ws.Range("H65").Value = 6902.75
$ws.Range("I65").Value = 6985
$ws.Range("J65").Value = 6656
$ws.Range("K65").Value = 34925
$ws.Range("L65").Value = 33280
$ws.Range("M65").Value = -31805
$ws.Range("N65").Value = -39520

# Row 122
$ws.Range("H122").Value = 1332.3889
$ws.Range("I122").Value = 1090.125
$ws.Range("J122").Value = 1526.2
$ws.Range("K122").Value = 3270.375
$ws.Range("L122").Value = 4578.6
$ws.Range("M122").Value = -820.375
$ws.Range("N122").Value = -9478.6

# Row 136
$ws.Range("H136").Value = 5390.9565
$ws.Range("I136").Value = 941.8823
$ws.Range("J136").Value = 17996.666
$ws.Range("K136").Value = 2825.6469
$ws.Range("L136").Value = 53989.99800000001
$ws.Range("M136").Value = -275.6468999999997

# Row 141
$ws.Range("H141").Value = 24360
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 24360
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 24360
$ws.Range("N141").Value = -34720

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 23037198
$ws.Range("I9").Value = 65013900
$ws.Range("J9").Value = 12543023
$ws.Range("K9").Value = 195041700
$ws.Range("L9").Value = 37629069
$ws.Range("M9").Value = -195041476
$ws.Range("N9").Value = -37629517

# Row 108
$ws.Range("H108").Value = 585.2
$ws.Range("I108").Value = 231.5
$ws.Range("J108").Value = 2000
$ws.Range("K108").Value = 694.5
$ws.Range("L108").Value = 6000
$ws.Range("M108").Value = 2185.5
$ws.Range("N108").Value = -11760

# Row 113
$ws.Range("H113").Value = 35108.45
$ws.Range("I113").Value = 966.6667
$ws.Range("J113").Value = 44015
$ws.Range("K113").Value = 2900.0001
$ws.Range("L113").Value = 132045
$ws.Range("M113").Value = -730.0001000000002

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4584.6665
$ws.Range("I70").Value = 4169.3335
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4169.3335
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -3899.3335

# Row 73
$ws.Range("H73").Value = 4584.6665
$ws.Range("I73").Value = 4169.3335
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4169.3335
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -3233.3335

# Row 122
$ws.Range("H122").Value = 6579947
$ws.Range("I122").Value = 13157894
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 39473682
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -39471232
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1102.6842
$ws.Range("I22").Value = 973.1539
$ws.Range("J22").Value = 1383.3334
$ws.Range("K22").Value = 973.1539
$ws.Range("L22").Value = 1383.3334
$ws.Range("M22").Value = -678.1539

# Row 27
$ws.Range("H27").Value = 1102.6842
$ws.Range("I27").Value = 973.1539
$ws.Range("J27").Value = 1383.3334
$ws.Range("K27").Value = 973.1539
$ws.Range("L27").Value = 1383.3334
$ws.Range("M27").Value = -866.1539

# Row 46
$ws.Range("H46").Value = 2049.8235
$ws.Range("I46").Value = 1766.5
$ws.Range("J46").Value = 2729.8
$ws.Range("K46").Value = 1766.5
$ws.Range("L46").Value = 2729.8
$ws.Range("M46").Value = -1578.5
$ws.Range("N46").Value = -3105.8

# Row 100
$ws.Range("H100").Value = 2128.5
$ws.Range("I100").Value = 2185.5715
$ws.Range("J100").Value = 2071.4285
$ws.Range("K100").Value = 2185.5715
$ws.Range("L100").Value = 2071.4285
$ws.Range("M100").Value = -1644.5715

# Row 119
$ws.Range("H119").Value = 40000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676

# Row 122
$ws.Range("H122").Value = 11289.909
$ws.Range("I122").Value = 14573.75
$ws.Range("J122").Value = 2533
$ws.Range("K122").Value = 43721.25
$ws.Range("L122").Value = 7599
$ws.Range("M122").Value = -41271.25
$ws.Range("N122").Value = -12499

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1575
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -10600

# Row 133
$ws.Range("H133").Value = 80405
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 80405
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 80405
$ws.Range("N133").Value = -90525
